$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$oldTotal = $wb.Worksheets.Item("总计")

# ------------------------------------------------------------------
# 1) Duplicate the existing "总计" sheet (placed right after it) so we
#    keep a full copy of the historical roll-up data/styles. This
#    duplicate will become the new "总计" sheet with the 2022-Q1 row
#    added on top.
# ------------------------------------------------------------------
$oldTotal.Copy([System.Reflection.Missing]::Value, $oldTotal)
$newTotal = $wb.Worksheets.Item("总计 (2)")
$newTotal.Name = "总计_newtmp"

# Insert a new data row right below the header and restore its
# formatting (ClearFormats removes the border/bold it inherits from
# row 1 on insert), then re-apply the style used by the other index
# cells in column A.
$newTotal.Rows("2:2").Insert(-4121)
$newTotal.Range("A2:D2").ClearFormats()
$newTotal.Range("A3").Copy()
$newTotal.Range("A2").PasteSpecial(-4122)

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 2
$newTotal.Range("D2").Value = 0.01

# Re-number the row index column and rewrite the value columns so
# that no residual floating point / shift artifacts remain.
$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q4"
$newTotal.Range("C3").Value = 2
$newTotal.Range("D3").Value = 0.04

$newTotal.Range("A4").Value = 2
$newTotal.Range("B4").Value = "2021-Q3"
$newTotal.Range("C4").Value = 7
$newTotal.Range("D4").Value = 0.76

$newTotal.Range("A5").Value = 3
$newTotal.Range("B5").Value = "2021-Q2"
$newTotal.Range("C5").Value = 3
$newTotal.Range("D5").Value = 0.19

# ------------------------------------------------------------------
# 2) Turn the original "总计" sheet into the new "2022-Q1" sheet (it
#    already sits in the right tab position, right after 2021-Q4 and
#    before 总计), replacing its contents with the fund holdings for
#    the new quarter.
# ------------------------------------------------------------------
$q1 = $oldTotal
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# Borrow the header styling (bold, bordered, centered) and the index
# column styling from the 2021-Q4 sheet, which uses the same 2-row
# layout.
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A3").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# The source data keeps these columns as plain text (matching the
# other quarterly sheets), so force a text format before assigning
# the values -- otherwise numeric-looking strings like "003854" or
# "0.19" would be auto-converted to numbers and lose their leading
# zeros / exact text representation.
$q1.Range("B2:G3").NumberFormat = "@"

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "003854"
$q1.Range("C2").Value = "汇安丰华灵活配置混合A"
$q1.Range("D2").Value = "0.19"
$q1.Range("E2").Value = "29.16"
$q1.Range("F2").Value = "2.02"
$q1.Range("G2").Value = "0.0038"
$q1.Range("H2").Value = 5

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "003855"
$q1.Range("C3").Value = "汇安丰华灵活配置混合C"
$q1.Range("D3").Value = "0.19"
$q1.Range("E3").Value = "29.16"
$q1.Range("F3").Value = "2.02"
$q1.Range("G3").Value = "0.0038"
$q1.Range("H3").Value = 5

# ------------------------------------------------------------------
# 3) Finally rename the duplicated roll-up sheet back to "总计".
# ------------------------------------------------------------------
$newTotal.Name = "总计"

Write-Host "done"
